$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Modification Type" value for the PKYVKQNTLKLAT row (row 2, col B)
# from the abbreviation "OX" to the full description "oxidized residue".
$ws.Range("B2").Value = "oxidized residue"

# Update the active selection shown when the workbook is saved.
$ws.Range("B5").Select()
